$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell content updates ---

# Row 2: Main Tool changed from Powerpoint to Photoshop
$ws.Range("F2").Value = "Photoshop"

# Row 2: Image URL updated to a new flickr link
$ws.Range("K2").Value = "http://farm9.staticflickr.com/8146/7469445394_ab04faaa6f.jpg"

# Row 4: new Image URL added (previously empty), match style used by other K column cells
$ws.Range("K4").Value = "http://farm8.staticflickr.com/7251/7469566482_a8e3a40df8.jpg"
$ws.Range("K4").WrapText = $true

# Row 4: Main Tool changed from SonyNex to SonyNex: Photoshop
$ws.Range("F4").Value = "SonyNex: Photoshop"

# Row 4: Phase Name changed from Construction to Development
$ws.Range("I4").Value = "Development"

# --- View/selection updates ---
# Scroll so column B is the leftmost visible column, and select F4 as the
# active cell (mirrors topLeftCell="B1" / selection activeCell="F4" sqref="F4").
$ws.Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("F4").Select()
